# Fixed iteration to next story through XPATH
#
# Adds two new worksheets at the end of the workbook:
#   - "garmin1": a duplicate of the existing "garmin" sheet (same 31 rows
#     of Garmin Ltd / Nikkei-225 story data).
#   - "garmin2": a fresh sheet with the header row plus one new story row
#     about Palo Alto Networks / Keysight / Garmin / Toll Brothers.

$wb = $excel.ActiveWorkbook

# --- garmin1: exact copy of the "garmin" sheet, placed after the last sheet ---
$srcSheet = $wb.Worksheets.Item("garmin")
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$srcSheet.Copy($null, $lastSheet)
$garmin1 = $wb.Worksheets.Item($wb.Worksheets.Count)
$garmin1.Name = "garmin1"

# --- garmin2: new sheet with header row + one new story row ---
$lastSheet2 = $wb.Worksheets.Item($wb.Worksheets.Count)
$garmin2 = $wb.Worksheets.Add($null, $lastSheet2)
$garmin2.Name = "garmin2"

$garmin2.Range("A1").Value = "Title"
$garmin2.Range("B1").Value = "Story_Content"
$garmin2.Range("C1").Value = "Date_Created"
$garmin2.Range("D1").Value = "Image_Filename"

$garmin2.Range("A2").Value = "Palo Alto Networks, Keysight fall; Garmin, Toll Brothers rise, Wednesday, 2/21/2024"
$garmin2.Range("B2").Value = "Stocks that traded heavily or had substantial price changes on Wednesday: Palo Alto Networks, Keysight fall; Garmin, Toll Brothers rise."
$garmin2.Range("C2").Value = "2024-02-21 21:24:35"
$garmin2.Range("D2").Value = "Palo_Alto_Networks_Keysight_fall;_Garmin_Toll_Brothers_rise_Wednesday_2/21/2024.png.png"

# Restore the originally active sheet/selection.
$wb.Worksheets.Item("Sheet").Activate()
